# corrección bisección 2.0 (lectura de potencias)
# Rewrites the Newton/bisection iteration table: the first three data rows
# (iterations 0-2) get corrected xn/fxn/E readings (the previous values were
# misread as complex numbers with "+0i" powers), and five more iteration
# rows (3-8) are appended so the table now spans A1:D10.
#
# All values in this table are stored as plain text (matching the workbook's
# inline-string cells), so every assignment is entered with a leading
# apostrophe to force Excel to keep the literal numeric-looking text instead
# of re-parsing it into a binary Number (which would both change the cell's
# type and silently round the displayed digits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iteration 0
$ws.Range("A2").Value = "'0"
$ws.Range("B2").Value = "'3.0"
$ws.Range("C2").Value = "'-65.0"
$ws.Range("D2").Value = "'1.000005"

# Iteration 1
$ws.Range("A3").Value = "'1"
$ws.Range("B3").Value = "'2.21424360312087"
$ws.Range("C3").Value = "'-26.6302712393383"
$ws.Range("D3").Value = "'0.354864476416074"

# Iteration 2
$ws.Range("A4").Value = "'2"
$ws.Range("B4").Value = "'1.16674765448471"
$ws.Range("C4").Value = "'-13.678949962338"
$ws.Range("D4").Value = "'0.897791347263332"

# Iteration 3 (new row)
$ws.Range("A5").Value = "'3"
$ws.Range("B5").Value = "'-1.77257584283229"
$ws.Range("C5").Value = "'-6.94363979939362"
$ws.Range("D5").Value = "'1.65822157015321"

# Iteration 4 (new row)
$ws.Range("A6").Value = "'4"
$ws.Range("B6").Value = "'-3.66772137943643"
$ws.Range("C6").Value = "'3.44598839162184"
$ws.Range("D6").Value = "'0.516709242754788"

# Iteration 5 (new row)
$ws.Range("A7").Value = "'5"
$ws.Range("B7").Value = "'-3.19849806177507"
$ws.Range("C7").Value = "'0.218523634218298"
$ws.Range("D7").Value = "'0.14670114178558"

# Iteration 6 (new row)
$ws.Range("A8").Value = "'6"
$ws.Range("B8").Value = "'-3.16442532963776"
$ws.Range("C8").Value = "'0.0011475026314019"
$ws.Range("D8").Value = "'0.0107674312356769"

# Iteration 7 (new row)
$ws.Range("A9").Value = "'7"
$ws.Range("B9").Value = "'-3.16424450940417"
$ws.Range("C9").Value = "'3.23050833219568e-08"
$ws.Range("D9").Value = "'5.71448360112314e-05"

# Iteration 8 (new row)
$ws.Range("A10").Value = "'8"
$ws.Range("B10").Value = "'-3.16424450431334"
$ws.Range("C10").Value = "'-1.00093544563862e-15"
$ws.Range("D10").Value = "'1.6088614078435902e-09"
